$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "angle" column header in A1, shift SIN header to B1, add D1 header ---
$ws.Range("B1").Value = "SIN"
$ws.Range("D1").Value = "Formula  Text"
$ws.Range("A1").Value = "angle"

# --- Column D: FORMULATEXT helper for rows 2-7 (showing the formula text used in column B) ---
$ws.Range("D2").Formula = "=FORMULATEXT(`$B2)"
$ws.Range("D3:D7").Formula = "=FORMULATEXT(`$B3)"

# --- Second example block: header row 9 (reuse the bold/centered header style from row 1) ---
$ws.Range("A1").Copy($ws.Range("A9"))
$ws.Range("A1").Copy($ws.Range("B9"))
$ws.Range("A1").Copy($ws.Range("C9"))
$ws.Range("A1").Copy($ws.Range("D9"))
$ws.Range("A9").Value = "Formula"
$ws.Range("B9").Value = "Formula Text"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "Comments"

# --- Error-propagation examples (rows 10-12) ---
$ws.Range("A10").Formula = "=SIN(SQRT(-1))"
$ws.Range("B10").Formula = "=FORMULATEXT(`$A10)"
$ws.Range("D10").Value = "Example of error propagation."

$ws.Range("A11").Formula = "=SIN(`"str`")"
$ws.Range("B11").Formula = "=FORMULATEXT(`$A11)"
$ws.Range("D11").Value = "Unable to convert angle argument to a number."

$ws.Range("A12").Formula = "=SIN(10/0)"
$ws.Range("B12").Formula = "=FORMULATEXT(`$A12)"
$ws.Range("D12").Value = "Input causes a #DIV/0! error."

# --- Selection / view state ---
$ws.Range("B13").Select()
